$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the extent of the data table.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

# Columns that may contain a HYPERLINK(...) formula that needs a friendly-name
# second argument added (S, T, U, V, W, X, Y => 19..25).
$hyperlinkCols = 19,20,21,22,23,24,25

for ($r = 2; $r -le $lastRow; $r++) {
    # Column C ("Förändrad") gets bumped from 45184 to 45186 for every data row.
    $cCell = $ws.Cells.Item($r, 3)
    if ($cCell.Value2 -eq 45184) {
        $cCell.Value2 = 45186
    }

    # The friendly name used in HYPERLINK(...) is the "Beteckning" value in column A.
    $name = $ws.Cells.Item($r, 1).Value2

    foreach ($col in $hyperlinkCols) {
        $cell = $ws.Cells.Item($r, $col)
        if ($cell.HasFormula) {
            $f = $cell.Formula
            if ($f -match '^=HYPERLINK\("([^"]*)"\)$') {
                $url = $matches[1]
                $cell.Formula = '=HYPERLINK("' + $url + '", "' + $name + '")'
            }
        }
    }
}
